# Add "Area" (column G) and "Atotal" (column H) to the discharge sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Per-segment incremental area, mirroring the Q (discharge) column pattern
# but without the velocity (C) factor.
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"

# Total cross-sectional area
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Match the author's final selection on the new total cell
$ws.Range("H2").Select() | Out-Null
